$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row (A1:L1) -------------------------------------------------
# Overwrite values in place (do NOT Clear the range) so the existing bold /
# centered / bordered header style already applied to row 1 is preserved.
$ws.Cells.Item(1, 1).Value = "Organization Name"
$ws.Cells.Item(1, 2).Value = "Categories"
$ws.Cells.Item(1, 3).Value = "Org URL"
$ws.Cells.Item(1, 4).Value = "Image URL"
$ws.Cells.Item(1, 5).Value = "Description"
$ws.Cells.Item(1, 6).Value = "Email"
$ws.Cells.Item(1, 7).Value = "Phone"
$ws.Cells.Item(1, 8).Value = "Website"
$ws.Cells.Item(1, 9).Value = "LinkedIn"
$ws.Cells.Item(1, 10).Value = "Instagram"
$ws.Cells.Item(1, 11).Value = "Facebook"
$ws.Cells.Item(1, 12).Value = "Twitter"

# --- 2. Data rows (A2:C70), columns D:L stay blank (already empty) --------
$ws.Cells.Item(2, 1).Value = "At a Glance"
$ws.Cells.Item(2, 2).Value = "General"
$ws.Cells.Item(2, 3).Value = "https://www.bscc.edu/about/at-a-glance"
$ws.Cells.Item(3, 1).Value = "Financial Data"
$ws.Cells.Item(3, 2).Value = "General"
$ws.Cells.Item(3, 3).Value = "https://www.bscc.edu/about/financial-data"
$ws.Cells.Item(4, 1).Value = "Hours of Operation"
$ws.Cells.Item(4, 2).Value = "General"
$ws.Cells.Item(4, 3).Value = "https://www.bscc.edu/about/hours-of-operation"
$ws.Cells.Item(5, 1).Value = "Employment Opportunities"
$ws.Cells.Item(5, 2).Value = "General"
$ws.Cells.Item(5, 3).Value = "https://www.bscc.edu/about/employment-opportunities"
$ws.Cells.Item(6, 1).Value = "Open Records Request"
$ws.Cells.Item(6, 2).Value = "General"
$ws.Cells.Item(6, 3).Value = "https://www.bscc.edu/about/open-records-request"
$ws.Cells.Item(7, 1).Value = "Economic Impact"
$ws.Cells.Item(7, 2).Value = "General"
$ws.Cells.Item(7, 3).Value = "https://www.bscc.edu/about/economic-impact"
$ws.Cells.Item(8, 1).Value = "Apply Online"
$ws.Cells.Item(8, 2).Value = "General"
$ws.Cells.Item(8, 3).Value = "https://www.bscc.edu/students/apply-online"
$ws.Cells.Item(9, 1).Value = "Current Students"
$ws.Cells.Item(9, 2).Value = "General"
$ws.Cells.Item(9, 3).Value = "https://www.bscc.edu/students/current-students"
$ws.Cells.Item(10, 1).Value = "Future Students"
$ws.Cells.Item(10, 2).Value = "General"
$ws.Cells.Item(10, 3).Value = "https://www.bscc.edu/students/future-students"
$ws.Cells.Item(11, 1).Value = "Transfer Students"
$ws.Cells.Item(11, 2).Value = "General"
$ws.Cells.Item(11, 3).Value = "https://www.bscc.edu/students/transfer-students"
$ws.Cells.Item(12, 1).Value = "Payment Plan"
$ws.Cells.Item(12, 2).Value = "General"
$ws.Cells.Item(12, 3).Value = "https://www.bscc.edu/students/payment-plan"
$ws.Cells.Item(13, 1).Value = "Scholarship Applications"
$ws.Cells.Item(13, 2).Value = "Academic"
$ws.Cells.Item(13, 3).Value = "https://www.bscc.edu/students/scholarship-applications"
$ws.Cells.Item(14, 1).Value = "Transcripts and Records"
$ws.Cells.Item(14, 2).Value = "General"
$ws.Cells.Item(14, 3).Value = "https://www.bscc.edu/students/transcripts-and-records"
$ws.Cells.Item(15, 1).Value = "International Students"
$ws.Cells.Item(15, 2).Value = "Cultural"
$ws.Cells.Item(15, 3).Value = "https://www.bscc.edu/students/international-students"
$ws.Cells.Item(16, 1).Value = "Student Resources"
$ws.Cells.Item(16, 2).Value = "General"
$ws.Cells.Item(16, 3).Value = "https://www.bscc.edu/students/student-resources"
$ws.Cells.Item(17, 1).Value = "Title IX"
$ws.Cells.Item(17, 2).Value = "General"
$ws.Cells.Item(17, 3).Value = "https://www.bscc.edu/students/title-ix"
$ws.Cells.Item(18, 1).Value = "Workforce Solutions"
$ws.Cells.Item(18, 2).Value = "General"
$ws.Cells.Item(18, 3).Value = "https://www.bscc.edu/programs/workforce-solutions"
$ws.Cells.Item(19, 1).Value = "Health Sciences"
$ws.Cells.Item(19, 2).Value = "General"
$ws.Cells.Item(19, 3).Value = "https://www.bscc.edu/programs/health-sciences"
$ws.Cells.Item(20, 1).Value = "Career Tech"
$ws.Cells.Item(20, 2).Value = "Professional"
$ws.Cells.Item(20, 3).Value = "https://www.bscc.edu/programs/career-tech"
$ws.Cells.Item(21, 1).Value = "Transfer Advisement"
$ws.Cells.Item(21, 2).Value = "General"
$ws.Cells.Item(21, 3).Value = "https://www.bscc.edu/programs/academic#transfer"
$ws.Cells.Item(22, 1).Value = "Alabama Transfers Guides"
$ws.Cells.Item(22, 2).Value = "General"
$ws.Cells.Item(22, 3).Value = "https://www.bscc.edu/students/current-students/student-organizations"
$ws.Cells.Item(23, 1).Value = "Advanced Engineering Design Technology (EDT)"
$ws.Cells.Item(23, 2).Value = "Professional"
$ws.Cells.Item(23, 3).Value = "https://www.bscc.edu/programs/career-tech/advanced-engineering-design-technology-edt"
$ws.Cells.Item(24, 1).Value = "Air Conditioning & Refrigeration"
$ws.Cells.Item(24, 2).Value = "General"
$ws.Cells.Item(24, 3).Value = "https://www.bscc.edu/programs/career-tech/air-conditioning-and-refrigeration"
$ws.Cells.Item(25, 1).Value = "Apprenticeship Program"
$ws.Cells.Item(25, 2).Value = "General"
$ws.Cells.Item(25, 3).Value = "https://www.bscc.edu/programs/career-tech/apprenticeship-program"
$ws.Cells.Item(26, 1).Value = "Automated Manufacturing Technology"
$ws.Cells.Item(26, 2).Value = "Professional"
$ws.Cells.Item(26, 3).Value = "https://www.bscc.edu/programs/career-tech/automated-manufacturing-technology"
$ws.Cells.Item(27, 1).Value = "Child Development"
$ws.Cells.Item(27, 2).Value = "General"
$ws.Cells.Item(27, 3).Value = "https://www.bscc.edu/programs/career-tech/child-development"
$ws.Cells.Item(28, 1).Value = "Computer Science"
$ws.Cells.Item(28, 2).Value = "Special Interest"
$ws.Cells.Item(28, 3).Value = "https://www.bscc.edu/programs/career-tech/computer-science"
$ws.Cells.Item(29, 1).Value = "Electrical Systems Technology"
$ws.Cells.Item(29, 2).Value = "Professional"
$ws.Cells.Item(29, 3).Value = "https://www.bscc.edu/programs/career-tech/electrical-systems-technology"
$ws.Cells.Item(30, 1).Value = "Industrial Electrical Technology (IST)"
$ws.Cells.Item(30, 2).Value = "Professional"
$ws.Cells.Item(30, 3).Value = "https://www.bscc.edu/programs/career-tech/industrial-electrical-technology"
$ws.Cells.Item(31, 1).Value = "Industrial Electronics Technology"
$ws.Cells.Item(31, 2).Value = "Professional"
$ws.Cells.Item(31, 3).Value = "https://www.bscc.edu/programs/career-tech/electrical-systems-technology"
$ws.Cells.Item(32, 1).Value = "Machine Tool Technology"
$ws.Cells.Item(32, 2).Value = "Professional"
$ws.Cells.Item(32, 3).Value = "https://www.bscc.edu/programs/career-tech/machine-tool-technology"
$ws.Cells.Item(33, 1).Value = "Management & Entrepreneurship"
$ws.Cells.Item(33, 2).Value = "General"
$ws.Cells.Item(33, 3).Value = "https://www.bscc.edu/programs/career-tech/management-and-entrepreneurship"
$ws.Cells.Item(34, 1).Value = "Mercedes Tech Program"
$ws.Cells.Item(34, 2).Value = "General"
$ws.Cells.Item(34, 3).Value = "https://www.bscc.edu/programs/career-tech/mercedes-tech-program"
$ws.Cells.Item(35, 1).Value = "Salon & Spa Management"
$ws.Cells.Item(35, 2).Value = "General"
$ws.Cells.Item(35, 3).Value = "https://www.bscc.edu/programs/career-tech/salon-spa-management"
$ws.Cells.Item(36, 1).Value = "Utility Vegetation Management"
$ws.Cells.Item(36, 2).Value = "General"
$ws.Cells.Item(36, 3).Value = "https://www.bscc.edu/programs/career-tech/utility-vegetation-management"
$ws.Cells.Item(37, 1).Value = "Vehicle Technology & Repair"
$ws.Cells.Item(37, 2).Value = "Professional"
$ws.Cells.Item(37, 3).Value = "https://www.bscc.edu/programs/career-tech/vehicle-technology-repair"
$ws.Cells.Item(38, 1).Value = "Welding Technology"
$ws.Cells.Item(38, 2).Value = "Professional"
$ws.Cells.Item(38, 3).Value = "https://www.bscc.edu/programs/career-tech/welding-technology"
$ws.Cells.Item(39, 1).Value = "Surgical Technology"
$ws.Cells.Item(39, 2).Value = "Professional"
$ws.Cells.Item(39, 3).Value = "https://www.bscc.edu/programs/health-sciences/surgical-technology"
$ws.Cells.Item(40, 1).Value = "Phlebotomy (10-Week Course)"
$ws.Cells.Item(40, 2).Value = "General"
$ws.Cells.Item(40, 3).Value = "https://www.bscc.edu/programs/health-sciences/phlebotomy-10-week-course"
$ws.Cells.Item(41, 1).Value = "Medical Assisting Technology"
$ws.Cells.Item(41, 2).Value = "Professional"
$ws.Cells.Item(41, 3).Value = "https://www.bscc.edu/programs/health-sciences/medical-assisting-technology"
$ws.Cells.Item(42, 1).Value = "Dental Assistant"
$ws.Cells.Item(42, 2).Value = "General"
$ws.Cells.Item(42, 3).Value = "https://www.bscc.edu/programs/health-sciences/dental-assistant"
$ws.Cells.Item(43, 1).Value = "BEST Robotics"
$ws.Cells.Item(43, 2).Value = "General"
$ws.Cells.Item(43, 3).Value = "https://www.bscc.edu/programs/workforce-solutions/best-robotics"
$ws.Cells.Item(44, 1).Value = "Manufacturing Skill Standards Council (MSSC)"
$ws.Cells.Item(44, 2).Value = "General"
$ws.Cells.Item(44, 3).Value = "http://www.bscc.edu/programs/workforce-solutions/mssc"
$ws.Cells.Item(45, 1).Value = "PeriOp 101"
$ws.Cells.Item(45, 2).Value = "General"
$ws.Cells.Item(45, 3).Value = "https://www.bscc.edu/programs/workforce-solutions/periop-101"
$ws.Cells.Item(46, 1).Value = "Second Opportunity System (SOS)"
$ws.Cells.Item(46, 2).Value = "General"
$ws.Cells.Item(46, 3).Value = "https://www.bscc.edu/programs/workforce-solutions/second-opportunity-system"
$ws.Cells.Item(47, 1).Value = "Final Exam Schedule"
$ws.Cells.Item(47, 2).Value = "General"
$ws.Cells.Item(47, 3).Value = "https://www.bscc.edu/students/current-students/final-exam-schedule"
$ws.Cells.Item(48, 1).Value = "Course Schedule"
$ws.Cells.Item(48, 2).Value = "General"
$ws.Cells.Item(48, 3).Value = "https://www.bscc.edu/students/current-students/course-schedule"
$ws.Cells.Item(49, 1).Value = "Canvas Assistance"
$ws.Cells.Item(49, 2).Value = "General"
$ws.Cells.Item(49, 3).Value = "https://www.bscc.edu/students/current-students/canvas-assistance"
$ws.Cells.Item(50, 1).Value = "Student ADA Accommodations"
$ws.Cells.Item(50, 2).Value = "General"
$ws.Cells.Item(50, 3).Value = "https://www.bscc.edu/students/current-students/student-ada-accommodations"
$ws.Cells.Item(51, 1).Value = "Course Withdrawal Form"
$ws.Cells.Item(51, 2).Value = "General"
$ws.Cells.Item(51, 3).Value = "https://www.bscc.edu/students/current-students/course-withdrawal-form"
$ws.Cells.Item(52, 1).Value = "Registration Assistance"
$ws.Cells.Item(52, 2).Value = "General"
$ws.Cells.Item(52, 3).Value = "https://www.bscc.edu/students/current-students/registration-assistance"
$ws.Cells.Item(53, 1).Value = "Forms, Publications, Applications"
$ws.Cells.Item(53, 2).Value = "General"
$ws.Cells.Item(53, 3).Value = "https://www.bscc.edu/students/current-students/forms-publications-applications"
$ws.Cells.Item(54, 1).Value = "Program of Study Change Form"
$ws.Cells.Item(54, 2).Value = "Academic"
$ws.Cells.Item(54, 3).Value = "https://www.bscc.edu/students/current-students/program-of-study-change-form"
$ws.Cells.Item(55, 1).Value = "Bear Alert"
$ws.Cells.Item(55, 2).Value = "General"
$ws.Cells.Item(55, 3).Value = "https://www.bscc.edu/students/current-students/bear-alert"
$ws.Cells.Item(56, 1).Value = "Free Speech Request Form"
$ws.Cells.Item(56, 2).Value = "General"
$ws.Cells.Item(56, 3).Value = "https://www.bscc.edu/students/current-students/free-speech-request-form"
$ws.Cells.Item(57, 1).Value = "Student Housing"
$ws.Cells.Item(57, 2).Value = "General"
$ws.Cells.Item(57, 3).Value = "https://www.bscc.edu/students/current-students/student-housing"
$ws.Cells.Item(58, 1).Value = "Student Organizations"
$ws.Cells.Item(58, 2).Value = "General"
$ws.Cells.Item(58, 3).Value = "https://www.bscc.edu/students/current-students/student-organizations"
$ws.Cells.Item(59, 1).Value = "MyBSCC Assistance"
$ws.Cells.Item(59, 2).Value = "General"
$ws.Cells.Item(59, 3).Value = "https://www.bscc.edu/students/current-students/mybscc-assistance"
$ws.Cells.Item(60, 1).Value = "ACCS Employment Opportunities"
$ws.Cells.Item(60, 2).Value = "General"
$ws.Cells.Item(60, 3).Value = "https://www.bscc.edu/students/current-students/student-organizations"
$ws.Cells.Item(61, 1).Value = "ACCS Formal Complaints (PDF)"
$ws.Cells.Item(61, 2).Value = "General"
$ws.Cells.Item(61, 3).Value = "https://www.bscc.edu/students/current-students/student-organizations"
$ws.Cells.Item(62, 1).Value = "ACCS Student Complaint Process"
$ws.Cells.Item(62, 2).Value = "General"
$ws.Cells.Item(62, 3).Value = "http://www.bscc.edu/accs-student-complaint-process"
$ws.Cells.Item(63, 1).Value = "Alabama Career Center (AlabamaWorks)"
$ws.Cells.Item(63, 2).Value = "Professional"
$ws.Cells.Item(63, 3).Value = "https://www.bscc.edu/students/current-students/student-organizations"
$ws.Cells.Item(64, 1).Value = "Bear Growl"
$ws.Cells.Item(64, 2).Value = "General"
$ws.Cells.Item(64, 3).Value = "https://www.bscc.edu/students/student-resources/technology"
$ws.Cells.Item(65, 1).Value = "Canvas|Canvas Assistance"
$ws.Cells.Item(65, 2).Value = "General"
$ws.Cells.Item(65, 3).Value = "https://www.bscc.edu/students/current-students/canvas-assistance"
$ws.Cells.Item(66, 1).Value = "Employee Webmail"
$ws.Cells.Item(66, 2).Value = "General"
$ws.Cells.Item(66, 3).Value = "https://www.bscc.edu/students/current-students/student-organizations"
$ws.Cells.Item(67, 1).Value = "Forms & Publications"
$ws.Cells.Item(67, 2).Value = "General"
$ws.Cells.Item(67, 3).Value = "https://www.bscc.edu/students/current-students/forms-publications-applications"
$ws.Cells.Item(68, 1).Value = "Live Chat"
$ws.Cells.Item(68, 2).Value = "General"
$ws.Cells.Item(68, 3).Value = "https://www.bscc.edu/live-chat"
$ws.Cells.Item(69, 1).Value = "Net Price Calculator"
$ws.Cells.Item(69, 2).Value = "General"
$ws.Cells.Item(69, 3).Value = "https://www.bscc.edu/students/current-students/student-organizations"
$ws.Cells.Item(70, 1).Value = "Operating Financial Data"
$ws.Cells.Item(70, 2).Value = "General"
$ws.Cells.Item(70, 3).Value = "https://www.bscc.edu/about/at-a-glance/financial-data"

# --- 3. Drop the old trailing rows 71:86 ------------------------------------
$ws.Range("A71:M86").ClearContents()

# --- 4. Drop column M (Tiktok Link) entirely --------------------------------
$ws.Columns.Item(13).Delete()

# --- 5. Column widths --------------------------------------------------------
# Excel's COM ColumnWidth is in characters; the stored OOXML width = 
# ColumnWidth + 0.8333333333333334 (5px padding / 6px-per-char for Calibri 11).
# Subtract that constant so the saved width matches the target exactly.
$ws.Columns.Item(1).ColumnWidth = 45.166666666666664
$ws.Columns.Item(2).ColumnWidth = 17.166666666666668
$ws.Columns.Item(7).ColumnWidth = 6.166666666666667
$ws.Columns.Item(8).ColumnWidth = 8.166666666666666
$ws.Columns.Item(9).ColumnWidth = 9.166666666666666
$ws.Columns.Item(10).ColumnWidth = 10.166666666666666
$ws.Columns.Item(11).ColumnWidth = 9.166666666666666
$ws.Columns.Item(12).ColumnWidth = 8.166666666666666
